{"js": "\nconst replacements = [\n  [\"735\u00d75=\", \"793\u00d75=\"],\n  [\"353\u00d73=\", \"849\u00d77=\"],\n  [\"276\u00d74=\", \"285\u00d74=\"],\n  [\"947\u00d78=\", \"817\u00d72=\"],\n  [\"217\u00d74=\", \"816\u00d73=\"],\n  [\"918\u00d76=\", \"905\u00d76=\"],\n  [\"913\u00d75=\", \"642\u00d73=\"],\n  [\"720\u00d74=\", \"619\u00d79=\"],\n  [\"660\u00d75=\", \"255\u00d76=\"],\n  [\"367\u00d72=\", \"550\u00d73=\"],\n  [\"998\u00d72=\", \"719\u00d76=\"],\n  [\"546\u00d74=\", \"108\u00d75=\"],\n  [\"561\u00d76=\", \"749\u00d77=\"],\n  [\"782\u00d76=\", \"527\u00d73=\"],\n  [\"778\u00d76=\", \"476\u00d78=\"],\n  [\"514\u00d75=\", \"626\u00d79=\"],\n  [\"898\u00d78=\", \"877\u00d78=\"],\n  [\"896\u00d76=\", \"395\u00d74=\"],\n  [\"910\u00d75=\", \"809\u00d75=\"],\n  [\"559\u00d78=\", \"690\u00d72=\"],\n  [\"762\u00d75=\", \"672\u00d74=\"],\n  [\"927\u00d75=\", \"888\u00d78=\"],\n  [\"212\u00d75=\", \"797\u00d77=\"],\n  [\"106\u00d78=\", \"614\u00d79=\"],\n  [\"333\u00d79=\", \"555\u00d75=\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\n\nreturn \"replaced:\" + totalReplaced;\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"735\u00d75=\", \"793\u00d75=\"),\n  @(\"353\u00d73=\", \"849\u00d77=\"),\n  @(\"276\u00d74=\", \"285\u00d74=\"),\n  @(\"947\u00d78=\", \"817\u00d72=\"),\n  @(\"217\u00d74=\", \"816\u00d73=\"),\n  @(\"918\u00d76=\", \"905\u00d76=\"),\n  @(\"913\u00d75=\", \"642\u00d73=\"),\n  @(\"720\u00d74=\", \"619\u00d79=\"),\n  @(\"660\u00d75=\", \"255\u00d76=\"),\n  @(\"367\u00d72=\", \"550\u00d73=\"),\n  @(\"998\u00d72=\", \"719\u00d76=\"),\n  @(\"546\u00d74=\", \"108\u00d75=\"),\n  @(\"561\u00d76=\", \"749\u00d77=\"),\n  @(\"782\u00d76=\", \"527\u00d73=\"),\n  @(\"778\u00d76=\", \"476\u00d78=\"),\n  @(\"514\u00d75=\", \"626\u00d79=\"),\n  @(\"898\u00d78=\", \"877\u00d78=\"),\n  @(\"896\u00d76=\", \"395\u00d74=\"),\n  @(\"910\u00d75=\", \"809\u00d75=\"),\n  @(\"559\u00d78=\", \"690\u00d72=\"),\n  @(\"762\u00d75=\", \"672\u00d74=\"),\n  @(\"927\u00d75=\", \"888\u00d78=\"),\n  @(\"212\u00d75=\", \"797\u00d77=\"),\n  @(\"106\u00d78=\", \"614\u00d79=\"),\n  @(\"333\u00d79=\", \"555\u00d75=\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
